$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2023" column (K), matching the formatting
# of the existing last data column (J).
$ws.Range("J3:J6").Copy($ws.Range("K3:K6"))

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 683
$ws.Range("K5").Value = 149
$ws.Range("K6").Value = 534

# The new column is now the right-hand edge of the table, so give it a
# thin right border (xlEdgeRight = 10) to close off the table.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

# Give the new column (and a few spare ones to its right) the same width
# as the other data columns.
$ws.Range("K1:O1").EntireColumn.ColumnWidth = 7.83

Write-Host "done"
